$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E3").Value = 12.7503
$ws.Range("A9").Value = -20.36849999999997
$ws.Range("A18").Value = -22.98810000000001
$ws.Range("A20").Value = -22.27040000000003
